$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12 & 13: the "At Work" (C) mark of 1.0 moves to "Sick Leave" (E)
$rowsCtoE = @(12, 13)
foreach ($r in $rowsCtoE) {
    $ws.Range("C$r").ClearContents()
    $ws.Range("E$r").Value = "'1.0"
}

# Rows 16-20, 23-27, 30-34: the "Sick Leave" (E) mark of 1.0 moves to "At Work" (C)
$rowsEtoC = @(16, 17, 18, 19, 20, 23, 24, 25, 26, 27, 30, 31, 32, 33, 34)
foreach ($r in $rowsEtoC) {
    $ws.Range("C$r").Value = "'1.0"
    $ws.Range("E$r").ClearContents()
}

# Totals row: reflect the new column totals after the re-tally
$ws.Range("C44").Value = "'18.0"
$ws.Range("E44").Value = "'2.0"
